$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13 mirrors the existing rows: copy row 12's formatting down to
# row 13 (so A13 picks up the same date style as A2:A12), then set values.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Range("A13").Value = 42619.891273148147
$ws.Range("B13").Value = 17
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = "Random"
